# Natmi LR-pair results were recomputed following Dr Hou's advice: a 4th
# "M2" sending/target cluster is now included, so the ECs/FAPs/sCs x
# ECs/FAPs/sCs grid of Plat-Lrp1 edges grows to a full 4x4 grid (rows 2-17)
# and every previously-computed statistic shifts to reflect the larger
# cell-count denominators.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Plat", "Lrp1", "ECs", 3, 1, 9.280752333333334, 27.842257, 0.1100677737582973, 0.1100677737582973, 3, 1, 8.131234, 24.393702, 0.02090995573015822, 0.02090995573015823, 75.46396891837934, 679.175720265414, 0.002301512276603067, 0.002301512276603068),
    @("ECs", "Plat", "Lrp1", "FAPs", 3, 1, 9.280752333333334, 27.842257, 0.1100677737582973, 0.1100677737582973, 3, 1, 243.3763986666667, 730.1291960000001, 0.625857000534647, 0.6258570005346471, 2258.71607980393, 20328.44471823538, 0.06888668673989407, 0.06888668673989408),
    @("ECs", "Plat", "Lrp1", "M2", 3, 1, 9.280752333333334, 27.842257, 0.1100677737582973, 0.1100677737582973, 3, 1, 103.9426383333333, 311.827915, 0.2672947262403034, 0.2672947262403035, 964.6658832449062, 8681.992949204157, 0.02942053545460372, 0.02942053545460373),
    @("ECs", "Plat", "Lrp1", "sCs", 3, 1, 9.280752333333334, 27.842257, 0.1100677737582973, 0.1100677737582973, 3, 1, 33.41874933333333, 100.256248, 0.08593831749489127, 0.08593831749489128, 310.1511358524151, 2791.360222671736, 0.009459039287196413, 0.009459039287196415),
    @("FAPs", "Plat", "Lrp1", "ECs", 3, 1, 66.457357, 199.372071, 0.7881702982826134, 0.7881702982826134, 3, 1, 8.131234, 24.393702, 0.02090995573015822, 0.02090995573015823, 540.3803207885379, 4863.422887096842, 0.01648060604491505, 0.01648060604491505),
    @("FAPs", "Plat", "Lrp1", "FAPs", 3, 1, 66.457357, 199.372071, 0.7881702982826134, 0.7881702982826134, 3, 1, 243.3763986666667, 730.1291960000001, 0.625857000534647, 0.6258570005346471, 16174.15221156499, 145567.3699040849, 0.4932818987936545, 0.4932818987936545),
    @("FAPs", "Plat", "Lrp1", "M2", 3, 1, 66.457357, 199.372071, 0.7881702982826134, 0.7881702982826134, 3, 1, 103.9426383333333, 311.827915, 0.2672947262403034, 0.2672947262403035, 6907.753023240219, 62169.77720916197, 0.2106737641101895, 0.2106737641101895),
    @("FAPs", "Plat", "Lrp1", "sCs", 3, 1, 66.457357, 199.372071, 0.7881702982826134, 0.7881702982826134, 3, 1, 33.41874933333333, 100.256248, 0.08593831749489127, 0.08593831749489128, 2220.921754938845, 19988.29579444961, 0.06773402933385439, 0.06773402933385439),
    @("M2", "Plat", "Lrp1", "ECs", 2, 0.6666666666666666, 0.3000473333333333, 0.900142, 0.003558498364782037, 0.003558498364782037, 3, 1, 8.131234, 24.393702, 0.02090995573015822, 0.02090995573015823, 2.439755078409333, 21.957795705684, 0.00007440804327343282, 0.00007440804327343283),
    @("M2", "Plat", "Lrp1", "FAPs", 2, 0.6666666666666666, 0.3000473333333333, 0.900142, 0.003558498364782037, 0.003558498364782037, 3, 1, 243.3763986666667, 730.1291960000001, 0.625857000534647, 0.6258570005346471, 73.02443941620356, 657.2199547458321, 0.002227111112989932, 0.002227111112989932),
    @("M2", "Plat", "Lrp1", "M2", 2, 0.6666666666666666, 0.3000473333333333, 0.900142, 0.003558498364782037, 0.003558498364782037, 3, 1, 103.9426383333333, 311.827915, 0.2672947262403034, 0.2672947262403035, 31.18771145154778, 280.68940306393, 0.000951167846240982, 0.0009511678462409823),
    @("M2", "Plat", "Lrp1", "sCs", 2, 0.6666666666666666, 0.3000473333333333, 0.900142, 0.003558498364782037, 0.003558498364782037, 3, 1, 33.41874933333333, 100.256248, 0.08593831749489127, 0.08593831749489128, 10.02720662080178, 90.244859587216, 0.0003058113622776901, 0.0003058113622776902),
    @("sCs", "Plat", "Lrp1", "ECs", 3, 1, 8.280368333333334, 24.841105, 0.09820342959430721, 0.09820342959430721, 3, 1, 8.131234, 24.393702, 0.02090995573015822, 0.02090995573015823, 67.32961252452333, 605.9665127207099, 0.002053429365366674, 0.002053429365366674),
    @("sCs", "Plat", "Lrp1", "FAPs", 3, 1, 8.280368333333334, 24.841105, 0.09820342959430721, 0.09820342959430721, 3, 1, 243.3763986666667, 730.1291960000001, 0.625857000534647, 0.6258570005346471, 2015.246224600176, 18137.21602140158, 0.0614613038881085, 0.06146130388810851),
    @("sCs", "Plat", "Lrp1", "M2", 3, 1, 8.280368333333334, 24.841105, 0.09820342959430721, 0.09820342959430721, 3, 1, 103.9426383333333, 311.827915, 0.2672947262403034, 0.2672947262403035, 860.6833309384529, 7746.149978446075, 0.02624925882926926, 0.02624925882926926),
    @("sCs", "Plat", "Lrp1", "sCs", 3, 1, 8.280368333333334, 24.841105, 0.09820342959430721, 0.09820342959430721, 3, 1, 33.41874933333333, 100.256248, 0.08593831749489127, 0.08593831749489128, 276.7195537193377, 2490.47598347404, 0.008439437511562774, 0.008439437511562776)
)

$numRows = $data.Count
$numCols = $data[0].Count
$arr = New-Object 'object[,]' $numRows, $numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $arr[$i, $j] = $data[$i][$j]
    }
}

# Write the refreshed 16-row result table (replaces + extends the prior 12 rows).
$ws.Range("A2:T17").Value = $arr

